$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "Games" sheet: append the now-completed game (was the first row of the
#    "Next" sheet, date 45306 vs WAS) as new row 41.
# ---------------------------------------------------------------------------
$games = $wb.Worksheets.Item("Games")

$newRow = 41
$games.Cells.Item($newRow, 1).Value = 40
$games.Cells.Item($newRow, 2).Value = 45306
$games.Cells.Item($newRow, 2).NumberFormat = $games.Range("B40").NumberFormat
$games.Cells.Item($newRow, 3).Value = 1
$games.Cells.Item($newRow, 4).Value = 129
$games.Cells.Item($newRow, 5).Value = 99.3
$games.Cells.Item($newRow, 6).Value = 0.622
$games.Cells.Item($newRow, 7).Value = 13
$games.Cells.Item($newRow, 8).Value = 32.4
$games.Cells.Item($newRow, 9).Value = 0.189
$games.Cells.Item($newRow, 10).Value = 129.9
$games.Cells.Item($newRow, 11).Value = "WAS"
$games.Cells.Item($newRow, 12).Value = 117
$games.Cells.Item($newRow, 13).Value = 0.5629999999999999
$games.Cells.Item($newRow, 14).Value = 7.4
$games.Cells.Item($newRow, 15).Value = 16.7
$games.Cells.Item($newRow, 16).Value = 0.105
$games.Cells.Item($newRow, 17).Value = 117.8
$games.Cells.Item($newRow, 18).Value = 0
$games.Cells.Item($newRow, 19).Value = 1

# ---------------------------------------------------------------------------
# 2) "Next" sheet: the game that was just played (row 2, WAS on 45306) comes
#    off the upcoming schedule. Delete that row and let everything below
#    shift up one; no new fixture is appended at the bottom.
# ---------------------------------------------------------------------------
$next = $wb.Worksheets.Item("Next")
$next.Rows.Item(2).Delete()
